$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = "Plastic"
